$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("P.J. Washington", "PF", "Dallas Mavericks"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Donovan Clingan", "C", "Portland Trail Blazers"),
    @("Jonathan Kuminga", "SF,PF", "Golden State Warriors"),
    @("Goga Bitadze", "C", "Orlando Magic"),
    @("Luguentz Dort", "SG,SF", "Oklahoma City Thunder")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
